$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Problema" row 5 gets a new "Comportamiento" note, and its old note
#     becomes the (longer) "Solucion" text, now spread over B5/C5 ---
$ws.Range("B5").Value = "No se alinean bien"
$ws.Range("C5").Value = "stylo de tabla definido en css.styles aunque no fue la solucion definitiva, tambien se declaro el estilo en el mismo html"

# --- row 3 ("Data frame muy grande/mostrar") gets a new comment + solution ---
$ws.Range("B3").Value = "Queda feo"
$ws.Range("C3").Value = "Lo puse abajo y listo"

# --- row 6 was the old "usuario ya existe" / Django forum entry -> delete it
#     entirely; this shifts every following row up by one (7->6, 11->10,
#     14->13, 17..25->16..24, 30->29), matching the target layout ---
$ws.Rows.Item(6).Delete()

# --- add the two brand-new problem/solution rows for the Boxplot topic ---
$ws.Range("A7").Value = "Cargan archivos no validos para analizar"
$ws.Range("B7").Value = "Se analiza igual tira cualquier cosa, el usuario puede pensar que los analisis no son fiables"
$ws.Range("C7").Value = "Leer si el archivo es csv"

$ws.Range("A8").Value = "Boxplot, no existe la variable"
$ws.Range("B8").Value = "Se rompe la web "
$ws.Range("C8").Value = "Error y mensaje"

# --- move the selection like the author left it ---
[void]$ws.Range("C12").Select()
